$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 29571
$ws.Range("J109").Value = 29571
$ws.Range("L109").Value = 29571
$ws.Range("N109").Value = -32345
$ws.Range("H112").Value = 1203.5682
$ws.Range("J112").Value = 1220.4048
$ws.Range("L112").Value = 3661.2144
$ws.Range("N112").Value = -5877.2144
$ws.Range("H114").Value = 40712
$ws.Range("J114").Value = 40712
$ws.Range("L114").Value = 40712
$ws.Range("N114").Value = -49390
$ws.Range("H116").Value = 9174.583000000001
$ws.Range("I116").Value = 4800
$ws.Range("K116").Value = 4800
$ws.Range("M116").Value = -1358
$ws.Range("H124").Value = 48517
$ws.Range("J124").Value = 48517
$ws.Range("L124").Value = 48517
$ws.Range("N124").Value = -58337
$ws.Range("H128").Value = 38000.8
$ws.Range("J128").Value = 38000.8
$ws.Range("L128").Value = 38000.8
$ws.Range("N128").Value = -47960.8
$ws.Range("H130").Value = 49087
$ws.Range("J130").Value = 49087
$ws.Range("L130").Value = 49087
$ws.Range("N130").Value = -59127
$ws.Range("H132").Value = 19127.736
$ws.Range("I132").Value = 2915.6743
$ws.Range("K132").Value = 8747.0229
$ws.Range("M132").Value = -6217.0229

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H111").Value = 47796.8
$ws.Range("J111").Value = 47796.8
$ws.Range("L111").Value = 47796.8
$ws.Range("N111").Value = -55976.8
$ws.Range("H113").Value = 33775
$ws.Range("J113").Value = 33775
$ws.Range("L113").Value = 33775
$ws.Range("N113").Value = -42453
$ws.Range("H114").Value = 36793.668
$ws.Range("J114").Value = 36793.668
$ws.Range("L114").Value = 36793.668
$ws.Range("N114").Value = -45471.668
$ws.Range("H121").Value = 39577.4
$ws.Range("J121").Value = 39577.4
$ws.Range("L121").Value = 39577.4
$ws.Range("N121").Value = -43071.4
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H125").Value = 50699
$ws.Range("J125").Value = 50699
$ws.Range("L125").Value = 50699
$ws.Range("N125").Value = -60539

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 546.4
$ws.Range("I94").Value = 389.63635
$ws.Range("J94").Value = 977.5
$ws.Range("K94").Value = 389.63635
$ws.Range("L94").Value = 977.5
$ws.Range("M94").Value = 61.36365000000001
$ws.Range("N94").Value = -1879.5
$ws.Range("H108").Value = 47684
$ws.Range("J108").Value = 47684
$ws.Range("L108").Value = 47684
$ws.Range("N108").Value = -55364
$ws.Range("H110").Value = 48084.668
$ws.Range("J110").Value = 48084.668
$ws.Range("L110").Value = 48084.668
$ws.Range("N110").Value = -56264.668
$ws.Range("H111").Value = 47694
$ws.Range("J111").Value = 47694
$ws.Range("L111").Value = 47694
$ws.Range("N111").Value = -55874
$ws.Range("H112").Value = 46487
$ws.Range("J112").Value = 46487
$ws.Range("L112").Value = 46487
$ws.Range("N112").Value = -49441
$ws.Range("H117").Value = 48935.5
$ws.Range("J117").Value = 48935.5
$ws.Range("L117").Value = 48935.5
$ws.Range("N117").Value = -58113.5
$ws.Range("H124").Value = 47996
$ws.Range("J124").Value = 47996
$ws.Range("L124").Value = 47996
$ws.Range("N124").Value = -57816
$ws.Range("H125").Value = 50570
$ws.Range("J125").Value = 50570
$ws.Range("L125").Value = 50570
$ws.Range("N125").Value = -60410
$ws.Range("H126").Value = 47181.332
$ws.Range("J126").Value = 47181.332
$ws.Range("L126").Value = 47181.332
$ws.Range("N126").Value = -57061.332
$ws.Range("H130").Value = 40552.855
$ws.Range("J130").Value = 40552.855
$ws.Range("L130").Value = 40552.855
$ws.Range("N130").Value = -50592.855

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H110").Value = 46816
$ws.Range("J110").Value = 46816
$ws.Range("L110").Value = 46816
$ws.Range("N110").Value = -54996
$ws.Range("H111").Value = 40189
$ws.Range("J111").Value = 40189
$ws.Range("L111").Value = 40189
$ws.Range("N111").Value = -48369
$ws.Range("H116").Value = 49366.5
$ws.Range("J116").Value = 49366.5
$ws.Range("L116").Value = 49366.5
$ws.Range("N116").Value = -58544.5
$ws.Range("H119").Value = 43501.332
$ws.Range("J119").Value = 43501.332
$ws.Range("L119").Value = 43501.332
$ws.Range("N119").Value = -53177.332

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 33706.4
$ws.Range("J110").Value = 33706.4
$ws.Range("L110").Value = 33706.4
$ws.Range("N110").Value = -41886.4
$ws.Range("H114").Value = 43134.6
$ws.Range("J114").Value = 43134.6
$ws.Range("L114").Value = 43134.6
$ws.Range("N114").Value = -51812.6
$ws.Range("H116").Value = 38996.8
$ws.Range("J116").Value = 38996.8
$ws.Range("L116").Value = 38996.8
$ws.Range("N116").Value = -48174.8
$ws.Range("H119").Value = 31049.334
$ws.Range("J119").Value = 31049.334
$ws.Range("L119").Value = 31049.334
$ws.Range("N119").Value = -40725.334
$ws.Range("H130").Value = 46374.727
$ws.Range("J130").Value = 46374.727
$ws.Range("L130").Value = 46374.727
$ws.Range("N130").Value = -56414.727

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 48626
$ws.Range("J108").Value = 48626
$ws.Range("L108").Value = 48626
$ws.Range("N108").Value = -56306
$ws.Range("H112").Value = 26609.2
$ws.Range("J112").Value = 29511.5
$ws.Range("L112").Value = 29511.5
$ws.Range("N112").Value = -32465.5
$ws.Range("H114").Value = 26263.334
$ws.Range("J114").Value = 26263.334
$ws.Range("L114").Value = 26263.334
$ws.Range("N114").Value = -34941.334
$ws.Range("H116").Value = 50676
$ws.Range("J116").Value = 50676
$ws.Range("L116").Value = 50676
$ws.Range("N116").Value = -59854
$ws.Range("H119").Value = 47412
$ws.Range("J119").Value = 47412
$ws.Range("L119").Value = 47412
$ws.Range("N119").Value = -57088
$ws.Range("H120").Value = 39212.75
$ws.Range("J120").Value = 39212.75
$ws.Range("L120").Value = 39212.75
$ws.Range("N120").Value = -48888.75
$ws.Range("H121").Value = 25558
$ws.Range("J121").Value = 25558
$ws.Range("L121").Value = 25558
$ws.Range("N121").Value = -29052
$ws.Range("H124").Value = 44714.5
$ws.Range("J124").Value = 44714.5
$ws.Range("L124").Value = 44714.5
$ws.Range("N124").Value = -54534.5
$ws.Range("H125").Value = 49707
$ws.Range("J125").Value = 49707
$ws.Range("L125").Value = 49707
$ws.Range("N125").Value = -59547
$ws.Range("H127").Value = 42915.75
$ws.Range("J127").Value = 42915.75
$ws.Range("L127").Value = 42915.75
$ws.Range("N127").Value = -52835.75
$ws.Range("H128").Value = 35210.5
$ws.Range("J128").Value = 35210.5
$ws.Range("L128").Value = 35210.5
$ws.Range("N128").Value = -45170.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 32864
$ws.Range("J108").Value = 32864
$ws.Range("L108").Value = 32864
$ws.Range("N108").Value = -40544
$ws.Range("H110").Value = 27632.5
$ws.Range("J110").Value = 27632.5
$ws.Range("L110").Value = 27632.5
$ws.Range("N110").Value = -35812.5
$ws.Range("H112").Value = 29377
$ws.Range("J112").Value = 29377
$ws.Range("L112").Value = 29377
$ws.Range("N112").Value = -32331
$ws.Range("H116").Value = 29583.5
$ws.Range("J116").Value = 29583.5
$ws.Range("L116").Value = 29583.5
$ws.Range("N116").Value = -38761.5
$ws.Range("H117").Value = 44704.5
$ws.Range("J117").Value = 44704.5
$ws.Range("L117").Value = 44704.5
$ws.Range("N117").Value = -53882.5
